$wb = $excel.ActiveWorkbook

# --- Sheet1: row height adjustments (rows grew taller due to wrapped text) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows.Item(8).RowHeight = 45
$ws1.Rows.Item(10).RowHeight = 30
$ws1.Rows.Item(11).RowHeight = 30
$ws1.Rows.Item(12).RowHeight = 30

# --- Sheet2: the reference/legend table (B4:I21) was reorganised and
# extended down to row 37 with several new entries. Clear the whole
# area first, then rewrite every cell to match the new layout. ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B4:I37").Clear()

$cells = @(
    @("B4", 25, 14),
    @("C4", "Loan TV Dispatch - Very", 14),
    @("I4", $null, 14),
    @("B5", 26, 14),
    @("C5", "Claim Closed - VERY", 14),
    @("I5", $null, 14),
    @("B6", 56, 14),
    @("C6", "Claim Closed - LW", 14),
    @("I6", $null, 14),
    @("E7", 41, 14),
    @("F7", "Courier Despatch - LW", 14),
    @("G7", $null, 14),
    @("H7", $null, 14),
    @("I7", $null, 14),
    @("B8", 27, 15),
    @("C8", "BER RPG - VERY", 15),
    @("B9", 31, 14),
    @("C9", "Failed Appt - Very SMS", 14),
    @("E9", 43, 14),
    @("F9", "Appointment Cancelled LW", 14),
    @("G9", $null, 14),
    @("H9", $null, 14),
    @("I9", $null, 14),
    @("B10", 32, 15),
    @("C10", "B2B Repair - VERY RPG", 15),
    @("E11", 42, 14),
    @("F11", "Courier Booked SMS LW", 14),
    @("G11", $null, 14),
    @("H11", $null, 14),
    @("I11", $null, 14),
    @("B12", 39, 14),
    @("C12", "Courier Booked LW", 14),
    @("E12", 45, 14),
    @("F12", "Failed Appointment - LW", 14),
    @("G12", $null, 14),
    @("H12", $null, 14),
    @("I12", $null, 14),
    @("B13", 48, 14),
    @("C13", "Claim Closed - LW RGP/MPI", 14),
    @("E13", 46, 14),
    @("F13", "Courier at Depot - LW", 14),
    @("G13", $null, 14),
    @("H13", $null, 14),
    @("I13", $null, 14),
    @("B14", 51, 14),
    @("C14", "Job Complete - LW", 14),
    @("E14", 47, 14),
    @("F14", "Courier at Depot - LW RPG", 14),
    @("G14", $null, 14),
    @("H14", $null, 14),
    @("I14", $null, 14),
    @("B16", 60, 14),
    @("C16", "BER", 14),
    @("E16", 50, 14),
    @("F16", "Delayed Appointment - LW", 14),
    @("G16", $null, 14),
    @("H16", $null, 14),
    @("I16", $null, 14),
    @("B17", 64, $null),
    @("C17", "B2BSMS", $null),
    @("E17", 59, 14),
    @("F17", "Delayed Appointment - SMS", 14),
    @("G17", $null, 14),
    @("H17", $null, 14),
    @("I17", $null, 14),
    @("B18", 65, $null),
    @("C18", "ApologyV", $null),
    @("B19", 66, $null),
    @("C19", "ApologyLW", $null),
    @("E19", 61, 14),
    @("F19", "Courier at Depot - SMS", 14),
    @("G19", $null, 14),
    @("H19", $null, 14),
    @("I19", $null, 14),
    @("B20", 68, $null),
    @("C20", "Delayed-SMS", $null),
    @("E20", 62, 14),
    @("F20", "Courier Despatched - SMS", 14),
    @("G20", $null, 14),
    @("H20", $null, 14),
    @("I20", $null, 14),
    @("B21", 73, $null),
    @("C21", "Annual health check very- extension", $null),
    @("E21", 63, 14),
    @("F21", "Courier Repair Complete- SMS", 14),
    @("G21", $null, 14),
    @("H21", $null, 14),
    @("I21", $null, 14),
    @("E23", 33, 14),
    @("F23", "New SG - Littlewoods", 14),
    @("G23", $null, 14),
    @("H23", $null, 14),
    @("I23", $null, 14),
    @("E25", 11, 14),
    @("F25", "Appointment Reminder Very", 14),
    @("G25", $null, 14),
    @("H25", $null, 14),
    @("I25", $null, 14),
    @("E26", 37, 14),
    @("F26", "Appointment Reminder Littlewoods", 14),
    @("G26", $null, 14),
    @("H26", $null, 14),
    @("I26", $null, 14),
    @("E28", 21, 14),
    @("F28", "Survey - Very", 14),
    @("G28", $null, 14),
    @("H28", $null, 14),
    @("I28", $null, 14),
    @("E29", 35, 14),
    @("F29", "Survey - Very RPG/MPI", 14),
    @("G29", $null, 14),
    @("H29", $null, 14),
    @("I29", $null, 14),
    @("E32", "Repair", $null),
    @("E33", 19, 14),
    @("F33", "Courier Repair Complete - Very", 14),
    @("G33", $null, 14),
    @("H33", $null, 14),
    @("I33", $null, 14),
    @("E34", 40, 14),
    @("F34", "Courier Repair Complete - LW", 14),
    @("G34", $null, 14),
    @("H34", $null, 14),
    @("I34", $null, 14),
    @("E35", 52, 14),
    @("F35", "B2B Repair - LW", 14),
    @("G35", $null, 14),
    @("H35", $null, 14),
    @("I35", $null, 14),
    @("E36", 54, 14),
    @("F36", "B2B Repair - LW RPG", 14),
    @("G36", $null, 14),
    @("H36", $null, 14),
    @("I36", $null, 14),
    @("E37", 18, 14),
    @("F37", "B2B Repair - Very", 14),
    @("G37", $null, 14),
    @("H37", $null, 14),
    @("I37", $null, 14)
)

$yellow = 65535
$red = 255

foreach ($item in $cells) {
    $ref = $item[0]
    $val = $item[1]
    $style = $item[2]
    $rng = $ws2.Range($ref)
    if ($null -ne $val) {
        $rng.Value = $val
    }
    if ($style -eq 14) {
        $rng.Interior.Color = $yellow
    } elseif ($style -eq 15) {
        $rng.Interior.Color = $red
    }
}
